$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price & volume/1h changes) per diff.
# Force text number format on target cells so numeric-looking strings
# (e.g. "266.80", "0.100") are preserved exactly as text, matching the
# original inline-string cell content instead of being coerced to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.782.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.288.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "119.45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.80"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.641"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.21%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.18"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0941"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.20"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.57"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.634.24"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.885"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.289.55"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.713.94"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.47"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.75"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.67"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.88"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.89"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "43.30"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.40"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.83"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.75"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0912"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.74"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0380"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.73"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.95"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.107"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.57"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.38"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.79"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.238"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.98"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.24%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.38"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.28"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.64"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.61"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +34.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.64"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.09%  "
